# Applies the cryptos-list value refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.968.20"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "2.753.88"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'575.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").Value = "'158.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.604"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.34%  "
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("E10").Value = "  +2.30%  "
$ws.Range("D11").Value = "'0.384"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("D12").Value = "'5.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -17.03%  "
$ws.Range("D13").Value = "3.240.98"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").Value = "'26.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("D15").Value = "63.584.93"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("D17").Value = "2.759.27"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "'12.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").Value = "'4.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.14%  "
$ws.Range("D20").Value = "'356.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("D21").Value = "'6.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.54%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'0.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").Value = "'0.532"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.04%  "
$ws.Range("D24").Value = "'65.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.54%  "
$ws.Range("D25").Value = "'0.170"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.53%  "
$ws.Range("D26").Value = "'8.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "0.0₃0905"
$ws.Range("E28").Value = "  -1.95%  "
$ws.Range("E29").Value = "  -2.69%  "
$ws.Range("D30").Value = "'7.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "'1.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("D32").Value = "'170.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("D33").Value = "'4.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("D34").Value = "'20.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.21%  "
$ws.Range("D35").Value = "'1.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").Value = "'1.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.58%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").Value = "'6.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("D40").Value = "'336.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").Value = "'4.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.50%  "
$ws.Range("D42").Value = "'39.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("D43").Value = "'21.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("D44").Value = "'21.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.10%  "
$ws.Range("D45").Value = "'0.0589"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("D46").Value = "'0.0255"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("D47").Value = "'0.102"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("E48").Value = "  -2.88%  "
$ws.Range("D49").Value = "'134.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("D50").Value = "'0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").Value = "'11.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.08%  "

Write-Output "Applied 97 cell updates"
